$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.772.97"
$ws.Range("E2").Value = "  +0.58%  "

$ws.Range("D3").Value = "1.644.25"
$ws.Range("E3").Value = "  +0.08%  "

$ws.Range("E4").Value = "  +0.46%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.70"
$ws.Range("E5").Value = "  +0.40%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.501"
$ws.Range("E6").Value = "  -0.43%  "

$ws.Range("E7").Value = "  +0.36%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.0630"
$ws.Range("E8").Value = "  +0.67%  "

$ws.Range("E9").Value = "  -0.19%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.17"
$ws.Range("E10").Value = "  -0.22%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0843"
$ws.Range("E11").Value = "  -0.07%  "

$ws.Range("D12").Value = "1.866.22"
$ws.Range("E12").Value = "  -0.34%  "

$ws.Range("D13").Value = "1.658.87"
$ws.Range("E13").Value = "  +0.95%  "

$ws.Range("E14").Value = "  -1.29%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.527"
$ws.Range("E15").Value = "  -0.85%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.47"
$ws.Range("E16").Value = "  -2.12%  "

$ws.Range("D17").Value = "26.746.86"
$ws.Range("E17").Value = "  +0.30%  "

$ws.Range("E18").Value = "  -1.63%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "213.79"
$ws.Range("E19").Value = "  -1.98%  "

$ws.Range("E20").Value = "  +0.35%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.37"
$ws.Range("E21").Value = "  +0.03%  "

$ws.Range("E22").Value = "  +12.93%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.26"
$ws.Range("E23").Value = "  -0.72%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.35"
$ws.Range("E24").Value = "  -2.06%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.07"
$ws.Range("E25").Value = "  -0.86%  "

$ws.Range("E26").Value = "  +0.36%  "

$ws.Range("E27").Value = "  -1.29%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.11"
$ws.Range("E28").Value = "  -0.06%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.66"
$ws.Range("E29").Value = "  -1.19%  "

$ws.Range("E31").Value = "  +0.32%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.32"
$ws.Range("E32").Value = "  -2.02%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.99"
$ws.Range("E33").Value = "  -2.04%  "

$ws.Range("D34").Value = "1.294.51"
$ws.Range("E34").Value = "  +1.44%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.54"
$ws.Range("E35").Value = "  -0.40%  "

$ws.Range("E36").Value = "  +1.45%  "

$ws.Range("E37").Value = "  -4.52%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.533"
$ws.Range("E38").Value = "  +1.10%  "

$ws.Range("E39").Value = "  -0.21%  "

$ws.Range("E40").Value = "  +0.40%  "

$ws.Range("E41").Value = "  +0.13%  "

$ws.Range("E42").Value = "  -0.04%  "

$ws.Range("E43").Value = "  -1.91%  "

$ws.Range("D44").Value = "1.793.83"
$ws.Range("E44").Value = "  +0.54%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.87"
$ws.Range("E45").Value = "  +3.59%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.47"
$ws.Range("E46").Value = "  -1.71%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.61"
$ws.Range("E47").Value = "  +0.05%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0524"
$ws.Range("E48").Value = "  +1.40%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.70"
$ws.Range("E49").Value = "  -1.12%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0975"
$ws.Range("E50").Value = "  -0.02%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.408"
$ws.Range("E51").Value = "  +0.11%  "
